$d = $word.ActiveDocument

# Each paragraph below is replaced in place (via Range.InsertXML) with the exact
# run / proofErr structure required by the edit, preserving the original paragraph
# identity attributes (w14:paraId/w14:textId/rsids) so only the runs change.

# Paragraph 1: "Oude KMSK Navbar: " -> split runs + proofErr around "Navbar"
$p1xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="1DA0171B" w14:textId="6F3CA8BD" w:rsidR="00EB6029" w:rsidRDefault="002939C1"><w:r><w:t xml:space="preserve">Oude KMSK </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Navbar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p>'
$d.Paragraphs(1).Range.InsertXML($p1xml) | Out-Null

# Paragraph 2: long "Home" navbar line -> split runs + proofErr around Snelschaak/k.schap/Elo
$p2xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6A01DE47" w14:textId="526CFF97" w:rsidR="002939C1" w:rsidRDefault="002939C1" w:rsidP="002939C1"><w:r><w:t xml:space="preserve">Home – Bestuur - Ligging lokaal – Kalender – Stapjestornooi – Jeugdkalender – Jeugdschaak – Clubkampioenschap – Laddertornooi – Interclub - Zilveren Toren - </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Snelschaak</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>k.schap</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> - Trofee Walter Huyck – Wintertornooi – Verzekering – </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Elo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – Historiek – Links - Inlog bestuur</w:t></w:r></w:p>'
$d.Paragraphs(2).Range.InsertXML($p2xml) | Out-Null

# Paragraph 4: "Nieuw KMSK Navbar:" -> split runs + proofErr around "Navbar"
$p4xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="41001558" w14:textId="621F76B6" w:rsidR="0068035C" w:rsidRDefault="0068035C" w:rsidP="002939C1"><w:r><w:t xml:space="preserve">Nieuw KMSK </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Navbar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p>'
$d.Paragraphs(4).Range.InsertXML($p4xml) | Out-Null

# Paragraph 5: "Nieuw KMSK" navbar detail line -> insert "Facebook" near the top,
#              replace the old inline "Facebook" mention with "Historiek", and
#              add the "Jeugdkalender" / "- Lesgevers" split runs
$p5xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3FE50700" w14:textId="4E18875C" w:rsidR="0068035C" w:rsidRDefault="0068035C" w:rsidP="002939C1"><w:r><w:t xml:space="preserve">Blog (Home) – </w:t></w:r><w:r><w:t xml:space="preserve"> Facebook - </w:t></w:r><w:r><w:t>Contact (Bestuur – Ligging Lokaal)</w:t></w:r><w:r><w:t xml:space="preserve"> – Toernooien (Stapjestornooi – Clubkampioenschap – Laddertoernooi – Snelschaaktoernooi</w:t></w:r><w:r><w:t xml:space="preserve"> – Trofee Walter Huyck</w:t></w:r><w:r><w:t xml:space="preserve"> - Wintertoernooi</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> – Team Competitie (Interclub – Zilveren Toren)</w:t></w:r><w:r><w:t xml:space="preserve"> – Links (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Elo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – Links</w:t></w:r><w:r><w:t xml:space="preserve"> -</w:t></w:r><w:r><w:t xml:space="preserve"> Historiek</w:t></w:r><w:r><w:t>) – Jeugdschaak (Inschrijven – Verzekering</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Jeugdkalender</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>- Lesgevers</w:t></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r><w:t>–</w:t></w:r><w:r><w:t xml:space="preserve"> Historiek</w:t></w:r></w:p>'
$d.Paragraphs(5).Range.InsertXML($p5xml) | Out-Null

# Paragraph 7: "HomePage: " -> proofErr around "HomePage"
$p7xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="605247E9" w14:textId="24692F88" w:rsidR="00A92E50" w:rsidRDefault="00A92E50" w:rsidP="002939C1"><w:proofErr w:type="spellStart"/><w:r><w:t>HomePage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p>'
$d.Paragraphs(7).Range.InsertXML($p7xml) | Out-Null

